$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells C11:C13 which held data before but are now superseded by column A
$ws.Range("C11:C13").ClearContents()

$cellData = @(
    @('A1', 'Cluster Name'),
    @('B1', 'Active cases'),
    @('C1', 'Cluster'),
    @('D1', 'Exist'),
    @('A2', '3612 BlueCross Glengowrie'),
    @('B2', 55),
    @('D2', 'old'),
    @('A3', '3684 Homestyle Aged Care Langford Grange Cranbourne East'),
    @('B3', 31),
    @('D3', 'old'),
    @('A4', '3980 Arcare Keysborough Aged Care Keysborough'),
    @('B4', 28),
    @('D4', 'old'),
    @('A5', '4518 Regis Aged Care Fawkner'),
    @('B5', 20),
    @('D5', 'old'),
    @('A6', 'Community Kids Pascoe Vale Early Education Centre Pascoe Vale'),
    @('B6', 25),
    @('D6', 'old'),
    @('A7', 'Guardian Childcare Caulfield'),
    @('B7', 19),
    @('D7', 'old'),
    @('A8', 'Inghams Enterprises Somerville'),
    @('B8', 29),
    @('D8', 'old'),
    @('A9', 'Oceania Meat Processors Laverton North'),
    @('B9', 15),
    @('D9', 'old'),
    @('A10', 'The Robin Hood Inn Drouin West'),
    @('B10', 48),
    @('D10', 'old'),
    @('A11', 'The Royal Children''s Hospital Melbourne Emergency Department Parkville Tier 1B'),
    @('B11', 12),
    @('D11', 'old'),
    @('A12', 'Werribee Mercy Hospital Emergency Department'),
    @('B12', 31),
    @('D12', 'old'),
    @('A13', 'Western Health Sunshine Hospital Emergency Department'),
    @('B13', 24),
    @('D13', 'old'),
    @('B14', 16),
    @('C14', '4075 Homestyle Aged Care Ferndale Gardens Aged Care Services Bayswater North'),
    @('D14', 'new'),
    @('B15', 22),
    @('C15', 'Western Health Sunshine Hospital Emergency Department'),
    @('D15', 'new'),
    @('B16', 23),
    @('C16', 'Community Kids Pascoe Vale Early Education Centre Pascoe Vale'),
    @('D16', 'new'),
    @('B17', 27),
    @('C17', '3980 Arcare Keysborough Aged Care Keysborough'),
    @('D17', 'new'),
    @('B18', 27),
    @('C18', 'Inghams Enterprises Somerville'),
    @('D18', 'new'),
    @('B19', 30),
    @('C19', '3684 Homestyle Aged Care Langford Grange Cranbourne East'),
    @('D19', 'new'),
    @('B20', 32),
    @('C20', '3824 Estia Health South Morang'),
    @('D20', 'new'),
    @('B21', 32),
    @('C21', 'St Vincents Hospital Emergency Department Melbourne'),
    @('D21', 'new'),
    @('B22', 49),
    @('C22', 'The Robin Hood Inn Drouin West'),
    @('D22', 'new'),
    @('B23', 56),
    @('C23', '3612 BlueCross Glengowrie'),
    @('D23', 'new')
)

foreach ($pair in $cellData) {
    $ws.Range($pair[0]).Value = $pair[1]
}
